$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.0118606375092661
$ws.Range("C2").Value = 0.00593031875463306
$ws.Range("D2").Value = 0.00593031875463306
$ws.Range("E2").Value = 0.00889547813194959
$ws.Range("F2").Value = 0.00148257968865827
$ws.Range("G2").Value = 0.011119347664937
$ws.Range("H2").Value = 0.116382505559674
$ws.Range("I2").Value = 0.0429948109710897
$ws.Range("J2").Value = 0.0148257968865827
$ws.Range("K2").Value = 0.0022238695329874
$ws.Range("L2").Value = 0.00889547813194959
$ws.Range("M2").Value = 0.00148257968865827
$ws.Range("N2").Value = 0.00148257968865827
$ws.Range("O2").Value = 0.000741289844329133
$ws.Range("P2").Value = 0.0185322461082283
$ws.Range("Q2").Value = 0.00148257968865827
$ws.Range("R2").Value = 0.00815418828762046
$ws.Range("S2").Value = 0.982209043736101
$ws.Range("T2").Value = 0.00815418828762046
$ws.Range("U2").Value = 0.0422535211267606
$ws.Range("V2").Value = 0.0289103039288362
$ws.Range("W2").Value = 0.00815418828762046
$ws.Range("X2").Value = 0.0022238695329874
$ws.Range("B3").Value = 0.00148257968865827
$ws.Range("C3").Value = 0.0022238695329874
$ws.Range("D3").Value = 0.687916975537435
$ws.Range("E3").Value = 0.0022238695329874
$ws.Range("F3").Value = 0.0378057820607858
$ws.Range("G3").Value = 0.00593031875463306
$ws.Range("H3").Value = 0.00518902891030393
$ws.Range("I3").Value = 0.916234247590808
$ws.Range("J3").Value = 0.97405485544848
$ws.Range("K3").Value = 0.00148257968865827
$ws.Range("L3").Value = 0.00518902891030393
$ws.Range("M3").Value = 0.00370644922164566
$ws.Range("N3").Value = 0.00518902891030393
$ws.Range("O3").Value = 0.983691623424759
$ws.Range("P3").Value = 0.00148257968865827
$ws.Range("Q3").Value = 0.000741289844329133
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.981467753891772
$ws.Range("U3").Value = 0.0896960711638251
$ws.Range("W3").Value = 0.0022238695329874
$ws.Range("X3").Value = 0.00296515937731653
$ws.Range("B4").Value = 0.98295033358043
$ws.Range("C4").Value = 0.99110452186805
$ws.Range("D4").Value = 0.00963676797627873
$ws.Range("E4").Value = 0.985174203113417
$ws.Range("F4").Value = 0.0192735359525575
$ws.Range("G4").Value = 0.979243884358784
$ws.Range("H4").Value = 0.865826538176427
$ws.Range("I4").Value = 0.00296515937731653
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.022238695329874
$ws.Range("L4").Value = 0.98295033358043
$ws.Range("M4").Value = 0.994069681245367
$ws.Range("N4").Value = 0.988880652335063
$ws.Range("O4").Value = 0.0022238695329874
$ws.Range("P4").Value = 0.972572275759822
$ws.Range("Q4").Value = 0.997776130467013
$ws.Range("R4").Value = 0.985915492957746
$ws.Range("S4").Value = 0.0170496664195701
$ws.Range("T4").Value = 0.00148257968865827
$ws.Range("U4").Value = 0.0259451445515196
$ws.Range("V4").Value = 0.966641957005189
$ws.Range("W4").Value = 0.987398072646405
$ws.Range("X4").Value = 0.992587101556709
$ws.Range("B5").Value = 0.00370644922164566
$ws.Range("C5").Value = 0.000741289844329133
$ws.Range("D5").Value = 0.296515937731653
$ws.Range("E5").Value = 0.00370644922164566
$ws.Range("F5").Value = 0.941438102297999
$ws.Range("G5").Value = 0.00370644922164566
$ws.Range("H5").Value = 0.0118606375092661
$ws.Range("I5").Value = 0.0370644922164566
$ws.Range("J5").Value = 0.011119347664937
$ws.Range("K5").Value = 0.97405485544848
$ws.Range("L5").Value = 0.0022238695329874
$ws.Range("M5").Value = 0.000741289844329133
$ws.Range("N5").Value = 0.00296515937731653
$ws.Range("O5").Value = 0.0133432171979244
$ws.Range("P5").Value = 0.00667160859896219
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.00593031875463306
$ws.Range("S5").Value = 0.000741289844329133
$ws.Range("T5").Value = 0.00889547813194959
$ws.Range("U5").Value = 0.842105263157895
$ws.Range("V5").Value = 0.0044477390659748
$ws.Range("W5").Value = 0.0022238695329874
$ws.Range("X5").Value = 0.00148257968865827
